# Updated cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - force text format to avoid numeric coercion
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D47", "D48", "D49", "D51")
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "98.732.94"
$ws.Range("D3").Value = "3.367.31"
$ws.Range("D5").Value = "258.70"
$ws.Range("D6").Value = "629.19"
$ws.Range("D8").Value = "0.390"
$ws.Range("D9").Value = "0.999"
$ws.Range("D10").Value = "0.859"
$ws.Range("D11").Value = "3.364.34"
$ws.Range("D13").Value = "98.323.77"
$ws.Range("D14").Value = "36.56"
$ws.Range("D15").Value = "0.0000248"
$ws.Range("D16").Value = "3.961.57"
$ws.Range("D17").Value = "5.49"
$ws.Range("D18").Value = "3.361.37"
$ws.Range("D19").Value = "3.56"
$ws.Range("D20").Value = "15.25"
$ws.Range("D21").Value = "489.41"
$ws.Range("D22").Value = "6.07"
$ws.Range("D23").Value = "0.0000210"
$ws.Range("D24").Value = "9.47"
$ws.Range("D25").Value = "5.65"
$ws.Range("D26").Value = "89.15"
$ws.Range("D27").Value = "11.92"
$ws.Range("D29").Value = "0.282"
$ws.Range("D30").Value = "1.00"
$ws.Range("D31").Value = "0.192"
$ws.Range("D32").Value = "0.135"
$ws.Range("D33").Value = "9.70"
$ws.Range("D34").Value = "0.999"
$ws.Range("D35").Value = "28.09"
$ws.Range("D36").Value = "0.151"
$ws.Range("D37").Value = "7.30"
$ws.Range("D39").Value = "499.51"
$ws.Range("D40").Value = "0.461"
$ws.Range("D47").Value = "159.49"
$ws.Range("D48").Value = "1.94"
$ws.Range("D49").Value = "0.851"
$ws.Range("D51").Value = "45.87"

foreach ($cell in $priceCells) {
    $ws.Range($cell).Style = "Normal"
}

# Volume(1h) (column E) updates
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("E3").Value = "  +7.43%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  +6.68%  "
$ws.Range("E6").Value = "  +3.06%  "
$ws.Range("E7").Value = "  +23.11%  "
$ws.Range("E8").Value = "  +1.86%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +10.38%  "
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("E14").Value = "  +7.76%  "
$ws.Range("E15").Value = "  +3.62%  "
$ws.Range("E16").Value = "  +6.70%  "
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("E18").Value = "  +7.44%  "
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("E20").Value = "  +5.52%  "
$ws.Range("E21").Value = "  -5.08%  "
$ws.Range("E22").Value = "  +7.73%  "
$ws.Range("E23").Value = "  +9.60%  "
$ws.Range("E24").Value = "  +8.18%  "
$ws.Range("E25").Value = "  +3.37%  "
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("E27").Value = "  +3.44%  "
$ws.Range("E28").Value = "  +7.78%  "
$ws.Range("E29").Value = "  +15.73%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("E31").Value = "  +9.55%  "
$ws.Range("E32").Value = "  +10.75%  "
$ws.Range("E33").Value = "  +8.30%  "
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("E35").Value = "  +5.67%  "
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("E38").Value = "  +4.73%  "
$ws.Range("E39").Value = "  +7.00%  "
$ws.Range("E40").Value = "  +6.42%  "
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("E42").Value = "  +4.42%  "
$ws.Range("E43").Value = "  +5.05%  "
$ws.Range("E44").Value = "  +5.51%  "
$ws.Range("E45").Value = "  +13.50%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("E49").Value = "  +13.20%  "
$ws.Range("E50").Value = "  +3.52%  "
$ws.Range("E51").Value = "  +3.70%  "
